$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: was a string "hh-10202" (SKU), now becomes a plain number (row count) 123
$ws.Range("B2").Value = 123

# Keep B3/B4 text the same (hh-10203 / hh-10204) - no content change required,
# but re-set them so the now-unused "hh-10202" shared string gets dropped.
$ws.Range("B3").Value = "hh-10203"
$ws.Range("B4").Value = "hh-10204"

# Update the active selection/view to C2
$ws.Range("C2").Select()
